{"js": "// Raise the reported test coverage numbers (79.07% -> 90.7%, etc.), rework\n// the surrounding sentences, split the paragraph so the int-overload\n// coverage gets its own sentence/paragraph, and append a closing remark\n// about the 90% coverage target to the narrative-summary paragraph.\n\nconst body = context.document.body;\n\n// --- Part 1: the \"two methods ... coverage of 79.07%...\" paragraph -------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet coverageParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"The two methods I chose for testing\") !== -1) {\n    coverageParaIndex = i;\n    break;\n  }\n}\nif (coverageParaIndex === -1) {\n  throw new Error(\"Could not find the 'two methods' paragraph\");\n}\n\nconst coverageParagraph = paragraphs.items[coverageParaIndex];\nconst newCoverageText =\n  \"The two methods I chose for testing were the getSquareArea(string) and getSquareArea(int) methods.  \" +\n  \"The string overloaded version only has a coverage of 90.7% due to only testing some of the string paths.  \" +\n  \"The 9.3% of the code not covered was from the other 4 string possibilities handled by if-statements.  \" +\n  \"I could have written tests for all ten, but I was shooting for the 90% code coverage. \";\n\ncoverageParagraph.insertText(newCoverageText, \"Replace\");\ncoverageParagraph.insertParagraph(\"The coverage for the int overload was 100%.\", \"After\");\nawait context.sync();\n\n// --- Part 2: the narrative-summary paragraph ending in \"...for testing?\" -\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet narrativeParaIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"is too much code for testing?\") !== -1) {\n    narrativeParaIndex = i;\n    break;\n  }\n}\nif (narrativeParaIndex === -1) {\n  throw new Error(\"Could not find the narrative-summary paragraph\");\n}\n\nconst narrativeParagraph = paragraphs2.items[narrativeParaIndex];\nnarrativeParagraph.load(\"text\");\nawait context.sync();\n\n// Rewriting the whole paragraph's text (rather than just appending a new\n// run) also clears the stale grammar-check <w:proofErr> markers that used\n// to bracket \"passed\", matching Word's own re-proofing behaviour.\nconst updatedNarrativeText =\n  narrativeParagraph.text +\n  \"  The 90% coverage seemed enough for this project, but I\\u2019m not convinced that it is for others.\";\n\nnarrativeParagraph.insertText(updatedNarrativeText, \"Replace\");\nawait context.sync();\n", "ps1": "# Raise the reported test coverage numbers (79.07% -> 90.7%, etc.), rework\n# the surrounding sentences, split the paragraph so the int-overload\n# coverage gets its own sentence/paragraph, and append a closing remark\n# about the 90% coverage target to the narrative-summary paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Part 1: the \"two methods ... coverage of 79.07%...\" paragraph -------\n$coverageParagraph = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*The two methods I chose for testing*\") {\n    $coverageParagraph = $p\n    break\n  }\n}\nif ($coverageParagraph -eq $null) {\n  throw \"Could not find the 'two methods' paragraph\"\n}\n\n$newCoverageText = \"The two methods I chose for testing were the getSquareArea(string) and getSquareArea(int) methods.  The string overloaded version only has a coverage of 90.7% due to only testing some of the string paths.  The 9.3% of the code not covered was from the other 4 string possibilities handled by if-statements.  I could have written tests for all ten, but I was shooting for the 90% code coverage. \"\n\n$cr = $coverageParagraph.Range\n$cr.MoveEnd(1, -1) | Out-Null\n$cr.Text = $newCoverageText\n\n# Insert a brand-new paragraph right after it for the int-overload coverage.\n$cr.InsertParagraphAfter()\n$intParagraph = $coverageParagraph.Next()\n$ir = $intParagraph.Range\n$ir.MoveEnd(1, -1) | Out-Null\n$ir.Text = \"The coverage for the int overload was 100%.\"\n\n# --- Part 2: the narrative-summary paragraph ending in \"...for testing?\" -\n$narrativeParagraph = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*is too much code for testing?*\") {\n    $narrativeParagraph = $p\n    break\n  }\n}\nif ($narrativeParagraph -eq $null) {\n  throw \"Could not find the narrative-summary paragraph\"\n}\n\n$nr = $narrativeParagraph.Range\n$nr.MoveEnd(1, -1) | Out-Null\n# Rewriting the whole paragraph's text (rather than just appending a new\n# run) also clears the stale grammar-check proofErr markers that used to\n# bracket \"passed\", matching Word's own re-proofing behaviour.\n$updatedNarrativeText = $nr.Text + \"  The 90% coverage seemed enough for this project, but I\u2019m not convinced that it is for others.\"\n$nr.Text = $updatedNarrativeText\n"}
